$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level permutation: the dataset was re-sorted/shuffled so that each
# destination row now holds the data that used to live in a different row.
# Columns A,B,C,E,F,G,H,I,J,K,R are constant across all data rows, so it is
# sufficient (and safe) to move the D:T span of each row according to the
# mapping below (destination row -> original source row).
$mapping = @{}
$mapping[2] = 21
$mapping[3] = 3
$mapping[4] = 4
$mapping[5] = 15
$mapping[6] = 16
$mapping[7] = 26
$mapping[8] = 27
$mapping[9] = 7
$mapping[10] = 48
$mapping[11] = 49
$mapping[12] = 34
$mapping[13] = 32
$mapping[14] = 23
$mapping[15] = 54
$mapping[16] = 42
$mapping[17] = 43
$mapping[18] = 6
$mapping[19] = 33
$mapping[20] = 9
$mapping[21] = 18
$mapping[22] = 19
$mapping[23] = 28
$mapping[24] = 37
$mapping[25] = 25
$mapping[26] = 8
$mapping[27] = 35
$mapping[28] = 22
$mapping[29] = 51
$mapping[30] = 59
$mapping[31] = 52
$mapping[32] = 53
$mapping[33] = 38
$mapping[34] = 39
$mapping[35] = 5
$mapping[36] = 56
$mapping[37] = 57
$mapping[38] = 44
$mapping[39] = 45
$mapping[40] = 11
$mapping[41] = 12
$mapping[42] = 30
$mapping[43] = 58
$mapping[44] = 17
$mapping[45] = 50
$mapping[46] = 20
$mapping[47] = 40
$mapping[48] = 41
$mapping[49] = 46
$mapping[50] = 47
$mapping[51] = 13
$mapping[52] = 14
$mapping[53] = 10
$mapping[54] = 36
$mapping[55] = 55
$mapping[56] = 24
$mapping[57] = 29
$mapping[58] = 31
$mapping[59] = 2

# Snapshot every row's D:T values BEFORE any writes, so overlapping moves
# (e.g. row 2 <- row 21, row 21 <- row 18, ...) do not clobber source data.
$snapshot = @{}
for ($r = 2; $r -le 59; $r++) {
    $snapshot[$r] = $ws.Range("D" + $r + ":T" + $r).Value2
}

# Write each destination row from its recorded source row.
for ($r = 2; $r -le 59; $r++) {
    $srcRow = $mapping[$r]
    $ws.Range("D" + $r + ":T" + $r).Value2 = $snapshot[$srcRow]
}
